$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = 0.55
$ws.Range("P2").Value = 0.1464285714285714
$ws.Range("S2").Value = 0.1035714285714286
$ws.Range("B3").Value = 0.006451612903225806
$ws.Range("C3").Value = 0.01935483870967742
$ws.Range("J3").Value = 0.006451612903225806
$ws.Range("P3").Value = 0.7741935483870968
$ws.Range("S3").Value = 0.1935483870967742
$ws.Range("J4").Value = 0.05882352941176471
$ws.Range("P4").Value = 0.7058823529411765
$ws.Range("S4").Value = 0.2352941176470588
$ws.Range("S5").Value = 1
$ws.Range("B6").Value = 0.05238095238095238
$ws.Range("D6").Value = 0.009523809523809525
$ws.Range("E6").Value = 0.004761904761904762
$ws.Range("F6").Value = 0.07142857142857142
$ws.Range("J6").Value = 0.2428571428571429
$ws.Range("O6").Value = 0.01904761904761905
$ws.Range("Q6").Value = 0.1809523809523809
$ws.Range("R6").Value = 0.0761904761904762
$ws.Range("S6").Value = 0.3428571428571429
$ws.Range("B7").Value = 0.1085714285714286
$ws.Range("D7").Value = 0.02285714285714286
$ws.Range("E7").Value = 0.005714285714285714
$ws.Range("F7").Value = 0.07428571428571429
$ws.Range("J7").Value = 0.09142857142857143
$ws.Range("O7").Value = 0.01714285714285714
$ws.Range("Q7").Value = 0.2514285714285714
$ws.Range("R7").Value = 0.03428571428571429
$ws.Range("S7").Value = 0.3942857142857143
$ws.Range("B8").Value = 0.1212121212121212
$ws.Range("D8").Value = 0.01948051948051948
$ws.Range("F8").Value = 0.05194805194805195
$ws.Range("J8").Value = 0.1017316017316017
$ws.Range("O8").Value = 0.01731601731601732
$ws.Range("Q8").Value = 0.1536796536796537
$ws.Range("R8").Value = 0.05844155844155844
$ws.Range("S8").Value = 0.4761904761904762
$ws.Range("B9").Value = 0.07964601769911504
$ws.Range("D9").Value = 0.01769911504424779
$ws.Range("F9").Value = 0.06637168141592921
$ws.Range("J9").Value = 0.05752212389380531
$ws.Range("O9").Value = 0.03097345132743363
$ws.Range("Q9").Value = 0.1769911504424779
$ws.Range("R9").Value = 0.06194690265486726
$ws.Range("S9").Value = 0.5088495575221239
$ws.Range("B10").Value = 0.1166150670794634
$ws.Range("D10").Value = 0.01754385964912281
$ws.Range("F10").Value = 0.08668730650154799
$ws.Range("J10").Value = 0.08771929824561403
$ws.Range("O10").Value = 0.01547987616099071
$ws.Range("Q10").Value = 0.1867905056759546
$ws.Range("R10").Value = 0.0804953560371517
$ws.Range("S10").Value = 0.4086687306501548
$ws.Range("G11").Value = 0.139622641509434
$ws.Range("J11").Value = 0.06415094339622641
$ws.Range("K11").Value = 0.2
$ws.Range("L11").Value = 0.569811320754717
$ws.Range("S11").Value = 0.02641509433962264
$ws.Range("G12").Value = 0.7770700636942676
$ws.Range("J12").Value = 0.1528662420382166
$ws.Range("K12").Value = 0.006369426751592357
$ws.Range("L12").Value = 0.02547770700636943
$ws.Range("S12").Value = 0.03821656050955414
$ws.Range("G13").Value = 0.7931034482758621
$ws.Range("J13").Value = 0.2068965517241379
$ws.Range("F15").Value = 0.05113636363636364
$ws.Range("H15").Value = 0.1534090909090909
$ws.Range("I15").Value = 0.1022727272727273
$ws.Range("J15").Value = 0.25
$ws.Range("K15").Value = 0.05681818181818182
$ws.Range("M15").Value = 0.01704545454545454
$ws.Range("O15").Value = 0.05681818181818182
$ws.Range("S15").Value = 0.3125
$ws.Range("F16").Value = 0.01136363636363636
$ws.Range("H16").Value = 0.1875
$ws.Range("I16").Value = 0.1079545454545455
$ws.Range("J16").Value = 0.4090909090909091
$ws.Range("K16").Value = 0.07954545454545454
$ws.Range("M16").Value = 0.01136363636363636
$ws.Range("O16").Value = 0.02840909090909091
$ws.Range("S16").Value = 0.1647727272727273
$ws.Range("F17").Value = 0.01333333333333333
$ws.Range("H17").Value = 0.216
$ws.Range("I17").Value = 0.1173333333333333
$ws.Range("J17").Value = 0.4106666666666667
$ws.Range("K17").Value = 0.06666666666666667
$ws.Range("M17").Value = 0.01866666666666667
$ws.Range("N17").Value = 0.002666666666666667
$ws.Range("O17").Value = 0.05866666666666667
$ws.Range("S17").Value = 0.096
$ws.Range("F18").Value = 0.02173913043478261
$ws.Range("H18").Value = 0.1884057971014493
$ws.Range("I18").Value = 0.108695652173913
$ws.Range("J18").Value = 0.3623188405797101
$ws.Range("K18").Value = 0.1159420289855072
$ws.Range("M18").Value = 0.02173913043478261
$ws.Range("O18").Value = 0.03623188405797102
$ws.Range("S18").Value = 0.1449275362318841
$ws.Range("F19").Value = 0.0129764801297648
$ws.Range("H19").Value = 0.2368207623682076
$ws.Range("I19").Value = 0.1062449310624493
$ws.Range("J19").Value = 0.326845093268451
$ws.Range("K19").Value = 0.1192214111922141
$ws.Range("M19").Value = 0.0129764801297648
$ws.Range("O19").Value = 0.0583941605839416
$ws.Range("S19").Value = 0.1265206812652068
